$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A1").Value = "LÍNEA 141 - LP1912"
$ws1.Range("A2").Value = "10/01 16:46:42"

# Columns: A=Hora_Scrap, B=Hora_Llegada, C=Linea, D=Minutos, E=Parada
$rows1 = @(
    @{ Row = 6;  A = "16:46:42"; B = "16:47"; C = "15_ABASTO";          D = 1 }
    @{ Row = 7;  A = "16:46:42"; B = "16:53"; C = "10_OLMOS";           D = 7 }
    @{ Row = 8;  A = "16:46:42"; B = "16:56"; C = "215C_EL PATO";       D = 10 }
    @{ Row = 9;  A = "16:46:42"; B = "17:01"; C = "16_SANTA ANA";       D = 15 }
    @{ Row = 10; A = "16:46:42"; B = "17:03"; C = "23_HERNANDEZ";       D = 17 }
    @{ Row = 11; A = "16:46:42"; B = "17:04"; C = "14_ABASTO";          D = 18 }
    @{ Row = 12; A = "16:46:42"; B = "17:07"; C = "15_ABASTO";          D = 21 }
    @{ Row = 13; A = "16:46:42"; B = "17:13"; C = "23_HERNANDEZ";       D = 27 }
    @{ Row = 14; A = "16:46:42"; B = "17:14"; C = "10_OLMOS";           D = 28 }
    @{ Row = 15; A = "16:46:42"; B = "17:17"; C = "17_ROMERO";          D = 31 }
    @{ Row = 16; A = "16:46:42"; B = "17:23"; C = "16_SANTA ANA";       D = 37 }
    @{ Row = 17; A = "16:46:42"; B = "17:24"; C = "11_ETCHEVERRY";      D = 38 }
    @{ Row = 18; A = "16:46:42"; B = "17:35"; C = "16_P MOR-SANTA ANA"; D = 49 }
    @{ Row = 19; A = "16:46:42"; B = "17:38"; C = "17X38_ROMERO";       D = 52 }
    @{ Row = 20; A = "16:46:42"; B = "17:44"; C = "215B_EL PATO";       D = 58 }
    @{ Row = 21; A = "16:46:42"; B = "17:48"; C = "27_EL RETIRO";       D = 62 }
    @{ Row = 22; A = "16:46:42"; B = "17:50"; C = "215_EL PELIGRO";     D = 64 }
    @{ Row = 23; A = "16:46:42"; B = "18:02"; C = "17_ROMERO";          D = 76 }
    @{ Row = 24; A = "16:46:42"; B = "18:04"; C = "14_ABASTO";          D = 78 }
    @{ Row = 25; A = "16:46:42"; B = "18:24"; C = "11_ETCHEVERRY";      D = 98 }
    @{ Row = 26; A = "16:46:42"; B = "18:34"; C = "14X44_ABASTO";       D = 108 }
    @{ Row = 27; A = "16:46:42"; B = "18:38"; C = "17X38_ROMERO";       D = 112 }
    @{ Row = 28; A = "16:46:42"; B = "18:41"; C = "16_P MOR-SANTA ANA"; D = 115 }
)

foreach ($r in $rows1) {
    $ws1.Cells.Item($r.Row, 1).Value = $r.A
    $ws1.Cells.Item($r.Row, 2).Value = $r.B
    $ws1.Cells.Item($r.Row, 3).Value = $r.C
    $ws1.Cells.Item($r.Row, 4).Value = $r.D
}

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A1").Value = "LÍNEA 141 - LP1912-215"
$ws2.Range("A2").Value = "10/01 16:46:42"

$rows2 = @(
    @{ Row = 6; A = "16:46:42"; B = "16:56"; C = "215C_EL PATO";   D = 10 }
    @{ Row = 7; A = "16:46:42"; B = "17:44"; C = "215B_EL PATO";   D = 58 }
    @{ Row = 8; A = "16:46:42"; B = "17:50"; C = "215_EL PELIGRO"; D = 64 }
)

foreach ($r in $rows2) {
    $ws2.Cells.Item($r.Row, 1).Value = $r.A
    $ws2.Cells.Item($r.Row, 2).Value = $r.B
    $ws2.Cells.Item($r.Row, 3).Value = $r.C
    $ws2.Cells.Item($r.Row, 4).Value = $r.D
}

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A1").Value = "LÍNEA 141 - 6203-6173"
$ws3.Range("A2").Value = "10/01 16:46:42"

$rows3 = @(
    @{ Row = 6; A = "16:46:42"; B = "16:59"; C = "215C_LA PLATA"; D = 13 }
    @{ Row = 7; A = "16:46:42"; B = "18:21"; C = "215C_LA PLATA"; D = 95 }
)

foreach ($r in $rows3) {
    $ws3.Cells.Item($r.Row, 1).Value = $r.A
    $ws3.Cells.Item($r.Row, 2).Value = $r.B
    $ws3.Cells.Item($r.Row, 3).Value = $r.C
    $ws3.Cells.Item($r.Row, 4).Value = $r.D
}
